# Update the "Generate Report for Handback" timestamps.
# Overview!G2 : "Latest HO Xliff Generate Date"      2016-09-07 06:16:51 -> 2016-09-07 06:18:47
# de-de!H2    : "Correspond Handoff Datetime"        2016-09-07 06:16:51 -> 2016-09-07 06:18:47  (shares value w/ Overview!G2)
# zh-cn!H2    : "Correspond Handoff Datetime"        2016-09-07 06:16:40 -> 2016-09-07 06:18:35
# zh-cn!K2    : "Correspond Handback DateTime"       2016-09-07 06:17:45 -> 2016-09-07 06:19:27
# de-de!K2    : "Correspond Handback DateTime"       2016-09-07 06:18:04 -> 2016-09-07 06:19:45

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-09-07 06:18:47"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-09-07 06:18:35"
$zhcn.Range("K2").Value = "2016-09-07 06:19:27"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-09-07 06:18:47"
$dede.Range("K2").Value = "2016-09-07 06:19:45"
